$wb = $excel.ActiveWorkbook

# --- Rename the "Requested quantity" header on the existing sheets ---
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Reuse the bold/bordered header style and the date number format already
# present in the workbook (copied from "Weekly Quantity") so the new sheet
# matches the existing look & feel instead of introducing brand new styles.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header row ---
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# --- Forecast data rows ---
$wsForecast.Range("A2").Value = 44934.99999999999
$wsForecast.Range("B2").Value = 10
$wsForecast.Range("C2").Value = 9.999999986452087
$wsForecast.Range("D2").Value = 10.00000001222594
$wsForecast.Range("A3").Value = 44941.99999999999
$wsForecast.Range("B3").Value = 10
$wsForecast.Range("C3").Value = 9.999999987043511
$wsForecast.Range("D3").Value = 10.00000001267165
$wsForecast.Range("A4").Value = 44948.99999999999
$wsForecast.Range("B4").Value = 10
$wsForecast.Range("C4").Value = 9.999999986316253
$wsForecast.Range("D4").Value = 10.00000001193312
$wsForecast.Range("A5").Value = 44955.99999999999
$wsForecast.Range("B5").Value = 10
$wsForecast.Range("C5").Value = 9.999999988397263
$wsForecast.Range("D5").Value = 10.00000001349635
$wsForecast.Range("A6").Value = 44962.99999999999
$wsForecast.Range("B6").Value = 10
$wsForecast.Range("C6").Value = 9.999999986907916
$wsForecast.Range("D6").Value = 10.00000001220072
$wsForecast.Range("A7").Value = 44969.99999999999
$wsForecast.Range("B7").Value = 10
$wsForecast.Range("C7").Value = 9.999999979407844
$wsForecast.Range("D7").Value = 10.00000001929287
$wsForecast.Range("A8").Value = 44976.99999999999
$wsForecast.Range("B8").Value = 10
$wsForecast.Range("C8").Value = 9.999999946089803
$wsForecast.Range("D8").Value = 10.00000004959863
$wsForecast.Range("A9").Value = 44983.99999999999
$wsForecast.Range("B9").Value = 10
$wsForecast.Range("C9").Value = 9.999999896004777
$wsForecast.Range("D9").Value = 10.00000009710475
$wsForecast.Range("A10").Value = 44990.99999999999
$wsForecast.Range("B10").Value = 10
$wsForecast.Range("C10").Value = 9.999999834396096
$wsForecast.Range("D10").Value = 10.00000015626566
$wsForecast.Range("A11").Value = 44997.99999999999
$wsForecast.Range("B11").Value = 10
$wsForecast.Range("C11").Value = 9.999999770174187
$wsForecast.Range("D11").Value = 10.00000022801089
$wsForecast.Range("A12").Value = 45004.99999999999
$wsForecast.Range("B12").Value = 10
$wsForecast.Range("C12").Value = 9.999999685601361
$wsForecast.Range("D12").Value = 10.00000029850147
$wsForecast.Range("A13").Value = 45011.99999999999
$wsForecast.Range("B13").Value = 10
$wsForecast.Range("C13").Value = 9.999999600159178
$wsForecast.Range("D13").Value = 10.0000003746063
$wsForecast.Range("A14").Value = 45018.99999999999
$wsForecast.Range("B14").Value = 10
$wsForecast.Range("C14").Value = 9.99999951975839
$wsForecast.Range("D14").Value = 10.00000047101607

Write-Output "PO Forecast sheet added and headers updated."
